$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header in C1
$ws.Range("C1").Value = "string_spouse"

# Fill in the new spouse data column C2:C10
$ws.Range("C2").Value = "Catelynken Thomasdr  Raes (1591 - ca. 1645)"
$ws.Range("C3").Value = 'c("Marya   Boscoop (? - 1642)", "Helena   Heussen (1621 - 1680)", "Maria   Nooseman (1652 - 1729)")'
$ws.Range("C4").Value = 'c("Geertje Gijsbertsdr, alias:  Giertje / Guertje Ghijsbers (1582 - ca. 1622)", "Aefje Willems  Saskers, alias:  Aafgie (1592 - ?)")'
$ws.Range("C5").Value = 'c("Duifke  de Bruyn (1638 - 1668)", "Niesje   Mangeles (1647 - 1710)")'
$ws.Range("C6").Value = 'c("Saertje Gerrits van Laar (1642 - 1683)", "Eva   Tol (1650 - 1727)")'
$ws.Range("C7").Value = "Paschasia  van Geldre (? - ?)"
$ws.Range("C8").Value = "Anneken Claesdr. de Moor (? - 1640)"
$ws.Range("C9").Value = 'c("Marguerite  van Bracht (? - 1600)", "Catharina  du Pire (ca. 1581 - 1654)")'
$ws.Range("C10").Value = 'c("Anneke Harmensdr.  Abeels (1590 - 1615)", "Lysbeth Reyniersdr. (1593 - 1675)")'
